# Add 4 new rows (293-296) of COVID overview data to the sheet,
# continuing the existing table that ends at row 292.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("2021-05-31", "overview", "K02000001", "United Kingdom", 4487339, 3383, 1, 127782),
    @("2021-06-01", "overview", "K02000001", "United Kingdom", 4490438, 3165, 0, 127782),
    @("2021-06-02", "overview", "K02000001", "United Kingdom", 4494699, 4330, 12, 127794),
    @("2021-06-03", "overview", "K02000001", "United Kingdom", 4499878, 5274, 18, 127812)
)

$startRow = 293
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Column A: date stored as plain text (matches existing rows, which use
    # inline strings rather than real date serials). Temporarily force a
    # text format so Excel doesn't auto-parse the string into a date
    # serial, then restore the default "Normal" style so the cell matches
    # its neighbours (no leftover custom number-format style).
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row[0]
    $cellA.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
